$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.178.95'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.828.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +10.41%  '
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '429.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.825.32'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.614'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.732'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.157'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000331'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.77'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +15.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.446.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +9.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +26.51%  '
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.818.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.46%  '
$ws.Range("E20").Value = '  +9.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '66.532.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '415.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +11.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.13'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +11.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '37.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +14.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +17.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +11.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +39.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.41%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '721.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.71%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +18.79%  '
$ws.Range("E33").Value = '  +15.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.77'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.37%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.83'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +47.98%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '39.22'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.150'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '55.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0470'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0731'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +19.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.89'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.71%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").Value = '  +7.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.137'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.52%  '
$ws.Range("E46").Value = '  +12.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.326'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +19.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +47.27%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.62%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.62%  '
$ws.Range("E51").Value = '  +7.90%  '
